$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of match data to append as row 3
$rowData = @(
    "z7ljIJn1",
    "14/11/2024",
    "16:30",
    "URUGUAY - PRIMERA DIVISION",
    "Boston River",
    "Rampla Juniors",
    1.73,
    3.6,
    4.75,
    2.38,
    2.1,
    5,
    1.06,
    10,
    1.36,
    3,
    2.1,
    1.7,
    1.44,
    2.63,
    1.91,
    1.8,
    6.5,
    7.5,
    8.5,
    13,
    15,
    29,
    9.5,
    7,
    17,
    51,
    12,
    23,
    15,
    51,
    41,
    41,
    351,
    3.6,
    9,
    21,
    29,
    51,
    151,
    2.63,
    8.5,
    67,
    6.5,
    26,
    34,
    101,
    126,
    301,
    51,
    51
)

$row = 3
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item($row, $col).Value = $rowData[$i]
}
